$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update header figures for the new account-statement period ---
$ws.Range("E11").Value = 260000
$ws.Range("C13").Value = 1

# --- Replace the worker block (rows 16-20) with the new worker's data:
#     CC 1104380660 - CARLOS ALFREDO RICARDO SOLA, periods 2410..2502 ---
$ws.Range("C16").Value = "1104380660"
$ws.Range("D16").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E16").Value = "2410"

$ws.Range("C17").Value = "1104380660"
$ws.Range("D17").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E17").Value = "2411"

$ws.Range("C18").Value = "1104380660"
$ws.Range("D18").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E18").Value = "2412"

$ws.Range("C19").Value = "1104380660"
$ws.Range("D19").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E19").Value = "2501"

$ws.Range("C20").Value = "1104380660"
$ws.Range("D20").Value = "CARLOS ALFREDO RICARDO SOLA"
$ws.Range("E20").Value = "2502"

# --- Remove the other two workers' blocks, previously rows 21-30
#     (MARIBEL BERMUDEZ DIAZ and a second CARLOS ALFREDO RICARDO SOLA
#     block) -- only one worker/period group remains now. ---
$ws.Range("A21:A30").EntireRow.Delete()

# --- Row 20 is now the last row of the (now shorter) data table, so it
#     gets the table's "closing" bottom border, same as before. ---
$ws.Range("B20:J20").Borders.Item(9).LineStyle = 1
$ws.Range("B20:J20").Borders.Item(9).Weight = 2
$ws.Range("B20:J20").Borders.Item(9).ColorIndex = 1
